# Update "想去人数" (want-to-go count) values in column F for the
# sheets that carry the full per-event listing: "展览" (sheet 1) and
# "全部类型" (sheet 4). The two sheets list the same events, but the
# "全部类型" sheet has one extra row inserted (a "演出" event) above
# row 16, so the row numbers differ by one from row 16 onward.

$wb = $excel.ActiveWorkbook

# row number (in "展览") -> new value for column F
$updates1 = @{
    3  = 1861
    4  = 37
    6  = 820
    14 = 137
    16 = 4385
    21 = 9
    22 = 995
    23 = 1876
    25 = 46
    26 = 18
    28 = 2075
    29 = 73
    31 = 6
    34 = 28
}

# row number (in "全部类型") -> new value for column F
$updates4 = @{
    3  = 1861
    4  = 37
    6  = 820
    14 = 137
    17 = 4385
    22 = 9
    23 = 995
    24 = 1876
    26 = 46
    27 = 18
    29 = 2075
    30 = 73
    32 = 6
    35 = 28
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
